$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.25"
$ws.Range("D4").Value = "'6.212"
$ws.Range("D5").Value = "'0.06187"
$ws.Range("D6").Value = "'3.580"
$ws.Range("D7").Value = "'1.531"
$ws.Range("D8").Value = "'6.549"
$ws.Range("D9").Value = "'0.8226"
$ws.Range("D10").Value = "'0.1646"
$ws.Range("D11").Value = "'0.08210"
$ws.Range("D12").Value = "'0.03466"
$ws.Range("D13").Value = "'0.03106"
$ws.Range("D14").Value = "'0.09127"
$ws.Range("D15").Value = "'3.774"
$ws.Range("D16").Value = "'0.001614"
$ws.Range("D17").Value = "'0.04695"
$ws.Range("D18").Value = "'0.006323"
$ws.Range("D19").Value = "'0.006134"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D22").Value = "'3.756"
$ws.Range("D23").Value = "'2.317"
$ws.Range("D24").Value = "'0.01384"
$ws.Range("D40").Value = "'0.04675"
$ws.Range("D41").Value = "'0.007006"
$ws.Range("D42").Value = "'0.1104"
$ws.Range("D43").Value = "'0.003522"
$ws.Range("D44").Value = "'0.01115"
$ws.Range("D45").Value = "'0.00006299"
$ws.Range("D47").Value = "'0.8457"
$ws.Range("D48").Value = "'0.001384"
